$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: round Q2/R2, clear Z2 and AB2 (time cells) ---
$ws.Range("Q2").Value2 = 333038
$ws.Range("R2").Value2 = 6626637
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# --- Row 3 & Row 4 swap their species/record identity data ---
# Save row 3 original values
$A3 = $ws.Range("A3").Value2
$B3 = $ws.Range("B3").Value2
$E3 = $ws.Range("E3").Value2
$F3 = $ws.Range("F3").Value2
$G3 = $ws.Range("G3").Value2
$H3 = $ws.Range("H3").Value2

# Save row 4 original values
$A4 = $ws.Range("A4").Value2
$B4 = $ws.Range("B4").Value2
$E4 = $ws.Range("E4").Value2
$F4 = $ws.Range("F4").Value2
$G4 = $ws.Range("G4").Value2
$H4 = $ws.Range("H4").Value2

# Write row 4's original data into row 3
$ws.Range("A3").Value2 = $A4
$ws.Range("B3").Value2 = $B4
$ws.Range("E3").Value2 = $E4
$ws.Range("F3").Value2 = $F4
$ws.Range("G3").Value2 = $G4
$ws.Range("H3").Value2 = $H4

# Write row 3's original data into row 4
$ws.Range("A4").Value2 = $A3
$ws.Range("B4").Value2 = $B3
$ws.Range("E4").Value2 = $E3
$ws.Range("F4").Value2 = $F3
$ws.Range("G4").Value2 = $G3
$ws.Range("H4").Value2 = $H3

# --- Row 3: round Q3/R3, clear Z3, move AB3 "00:00" -> AC3 comment ---
$ws.Range("Q3").Value2 = 333038
$ws.Range("R3").Value2 = 6626637
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
$ws.Range("AC3").Value2 = "Örtrikt dråg i granskog"

# --- Row 4: round Q4/R4, clear Z4, AB4, AC4 ---
$ws.Range("Q4").Value2 = 333038
$ws.Range("R4").Value2 = 6626637
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
$ws.Range("AC4").ClearContents()
